$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "nome" column (B) -- it was redundant with "cognome" now that
# the sheet only needs the surname + arrival status.
$ws.Columns("B").Delete()

# Remove the "nome completo" helper column (now shifted into column C)
# -- it was just a concatenation of the two name columns and is no longer
# needed.
$ws.Columns("C").Delete()

# Fix the demo data: correct which people have actually arrived.
$ws.Range("C2").Value = "Arrivato"
$ws.Range("C3").Value = "No"
$ws.Range("C8").Value = "No"
